$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value edits -----------------------------------------------------
$ws.Range("D1").Value = "Tandems"

$ws.Range("G2").Value = 140
$ws.Range("G3").Value = 140
$ws.Range("G4").Value = 150
$ws.Range("G5").Value = 140
$ws.Range("G6").Value = 150

$ws.Range("H3").Value = "Mizumoto et al., 2021  Proc R Soc B, 288, 20210998"
$ws.Range("H5").Value = "Mizumoto et al., 2021  Proc R Soc B, 288, 20210998"

# --- Formatting edits --------------------------------------------------
$used = $ws.Range("A1:H6")
$used.Font.Name = "PT Serif"
$used.Font.Size = 9
$used.Interior.ThemeColor = 2

Write-Host "done"
